$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths (engine quantizes ColumnWidth to nearest 1/6 of a character,
# so we dial in the input that lands closest to the desired stored width).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth  = 14.592447916666666   # A -> 15.42578125
$ws.Columns.Item(2).ColumnWidth  = 18.166666666666668   # B -> 19
$ws.Columns.Item(4).ColumnWidth  = 13.307291666666666   # D -> 14.140625
$ws.Columns.Item(5).ColumnWidth  = 10.307291666666666   # E -> 11.140625
$ws.Columns.Item(7).ColumnWidth  = 10.451822916666666   # G -> 11.28515625
$ws.Columns.Item(8).ColumnWidth  = 6.022135416666667    # H -> 6.85546875
$ws.Columns.Item(9).ColumnWidth  = 16.736979166666668   # I -> 17.5703125
$ws.Columns.Item(10).ColumnWidth = 18.736979166666668   # J -> 19.5703125
$ws.Columns.Item(12).ColumnWidth = 13.592447916666666   # L -> 14.42578125
$ws.Columns.Item(13).ColumnWidth = 36.592447916666664   # M -> 37.42578125

# ---------------------------------------------------------------------------
# Row 1 - headers
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Название"
$ws.Range("B1").Value = "Жанры"
$ws.Range("C1").Value = "Теги"
$ws.Range("D1").Value = "Магазин"
$ws.Range("E1").Value = "Стоимость"
$ws.Range("F1").Value = "ОС"
$ws.Range("G1").Value = "Процессор"
$ws.Range("H1").Value = "ОЗУ"
$ws.Range("I1").Value = "Видеокарта"
$ws.Range("J1").Value = "Память видеокарты"
$ws.Range("K1").Value = "DirectX"
$ws.Range("L1").Value = "Жесткий диск"

$ws.Range("A1:L1").Font.Bold = $true

# M1 is a rich-text run: bold "Прочие" + space + bold "требования"
$ws.Range("M1").Value = "Прочие требования"
$m1 = $ws.Range("M1")
[void]$m1.Characters(1, 6).Font.Bold
$m1.Characters(1, 6).Font.Bold = $true
$m1.Characters(1, 6).Font.Name = "Calibri"
$m1.Characters(1, 6).Font.Size = 11
$m1.Characters(8, 10).Font.Bold = $true
$m1.Characters(8, 10).Font.Name = "Calibri"
$m1.Characters(8, 10).Font.Size = 11

# ---------------------------------------------------------------------------
# Row 2 - Apex Legends
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Apex Legends"
$ws.Range("B2").Value = "Шутер / Королевская битва"
$ws.Range("D2").Value = "Origin и Steam"
$ws.Range("E2").Value = "Free to Play"
$ws.Range("F2").Value = "64-разрядная версия Windiws 7"
$ws.Range("G2").Value = "Четырехъядерный процессор Intel Core i3-6300 3,8 ГГц "
$ws.Range("H2").Value = "6 ГБ"
$ws.Range("I2").Value = "nVidia GeForce GT 640 / Radeon HD 7730"
$ws.Range("J2").Value = "1 ГБ"
$ws.Range("L2").Value = "не менее 22 ГБ"

# ---------------------------------------------------------------------------
# Row 3 - Star Wars Battlefront 2
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Star Wars Battlefront 2"
$ws.Range("B3").Value = "Экшен / Шутер"
$ws.Range("D3").Value = "Steam"
$ws.Range("E3").Value = "999 руб (Стандартное издание)"
$ws.Range("F3").Value = "64-разрядная Windows 7 SP1, Windows 8.1 или Windows 10"
$ws.Range("G3").Value = "Intel Core i5 6600K / AMD FX-6350"
$ws.Range("H3").Value = "8 ГБ"
$ws.Range("I3").Value = "nVidia GeForce GTX 660 / AMD Radeon HD 7850"
$ws.Range("J3").Value = "2 ГБ"
$ws.Range("K3").Value = "11.0 или аналогичная"
$ws.Range("L3").Value = "не менее 55 ГБ"
$ws.Range("M3").Value = "Требования к соединению с интернетом: скорость 512 кбит/с"

# ---------------------------------------------------------------------------
# Row 4 - Minecraft Windows 10 edition
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Minecraft Windows 10 edition"
$ws.Range("B4").Value = "Экшн / Приключения"
$ws.Range("D4").Value = "Microsoft Store"
$ws.Range("E4").Value = "29.99 $"
$ws.Range("F4").Value = "Windows Xp 32"
$ws.Range("G4").Value = "Intel Pentium 4 1,6 Ггц / AMD Athlon XP 1600+"
$ws.Range("H4").Value = "500 МБ"
$ws.Range("I4").Value = "nVidia GeForce 4 MX 440 / AMD Radeon HD 3200"
$ws.Range("K4").Value = "8 или аналогичная"
$ws.Range("L4").Value = "не менее 1 ГБ"

# Row 6 stays empty but keeps a custom (slightly reduced) row height, as in
# the source workbook.
$ws.Rows.Item(6).RowHeight = 13.5

# Selection / view state
[void]$ws.Range("L4").Select()

Write-Host "edit complete"
